$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "245.92"
    "D3" = "22.28"
    "D4" = "5.350"
    "D5" = "0.05916"
    "D7" = "6.390"
    "D8" = "0.8122"
    "D9" = "0.9602"
    "D10" = "0.1429"
    "D11" = "0.07407"
    "D13" = "0.03042"
    "D14" = "4.452"
    "D15" = "0.09401"
    "D16" = "0.001597"
    "D17" = "0.04829"
    "D18" = "0.0005911"
    "D19" = "0.006133"
    "D20" = "0.004087"
    "D21" = "0.0009834"
    "D22" = "0.00009701"
    "D23" = "3.727"
    "D24" = "2.164"
    "D26" = "0.1333"
    "D27" = "0.0002461"
    "D40" = "0.03941"
    "D41" = "0.1074"
    "D42" = "0.002430"
    "D43" = "0.003046"
    "D44" = "0.005344"
    "D45" = "0.00005302"
    "D47" = "0.7501"
    "D48" = "0.04771"
    "D49" = "0.00002100"
}

foreach ($addr in $updates.Keys) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $updates[$addr]
}
